# Generate Report for Handback
# Update the timestamp values on the report sheets (these are stored as
# plain text strings, formatted like "yyyy-MM-dd HH:mm:ss").

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), row 2
$wsOverview.Range("G2").Value = "2016-09-06 05:06:52"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns, row 2
$wsZhCn.Range("H2").Value = "2016-09-06 05:06:38"
$wsZhCn.Range("K2").Value = "2016-09-06 05:07:46"

# de-de sheet: "Correspond Handback DateTime" (K) column, row 2
$wsDeDe.Range("K2").Value = "2016-09-06 05:08:08"
